$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells that need updated numeric-looking text values
# so Excel does not auto-convert them into numbers (matching the original t="inlineStr" text cells).
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.514.60"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.63"
$ws.Range("E3").Value = "  +0.56%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.82%  "

# Row 5
$ws.Range("E5").Value = "  -0.19%  "

# Row 6
$ws.Range("E6").Value = "  +0.75%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4808"
$ws.Range("E7").Value = "  -0.68%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4052"
$ws.Range("E8").Value = "  -0.45%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08222"
$ws.Range("E9").Value = "  +1.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.009"
$ws.Range("E10").Value = "  -0.03%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.47"
$ws.Range("E11").Value = "  +0.18%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.901.14"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.053"
$ws.Range("E13").Value = "  +0.55%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.236"
$ws.Range("E14").Value = "  +2.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.51"
$ws.Range("E15").Value = "  +1.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06855"
$ws.Range("E16").Value = "  +1.38%  "

# Row 17
$ws.Range("E17").Value = "  +0.74%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001038"
$ws.Range("E18").Value = "  -0.30%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.50"
$ws.Range("E19").Value = "  -1.06%  "

# Row 20
$ws.Range("E20").Value = "  +0.67%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.516.27"
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.670"
$ws.Range("E22").Value = "  +1.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.86"
$ws.Range("E23").Value = "  +0.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.192"
$ws.Range("E24").Value = "  +1.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.139.57"
$ws.Range("E25").Value = "  +0.50%  "

# Row 26
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.532"
$ws.Range("E26").Value = "  +3.89%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.98"
$ws.Range("E27").Value = "  +0.59%  "

# Row 28
$ws.Range("E28").Value = "  -0.02%  "

# Row 29
$ws.Range("E29").Value = "  -0.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.49"
$ws.Range("E30").Value = "  +1.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.013"
$ws.Range("E31").Value = "  -2.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09621"
$ws.Range("E32").Value = "  +0.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.612"
$ws.Range("E33").Value = "  +1.23%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.560"
$ws.Range("E34").Value = "  +0.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.371"
$ws.Range("E35").Value = "  -1.52%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06312"
$ws.Range("E36").Value = "  +3.32%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02278"
$ws.Range("E37").Value = "  +0.62%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.181"
$ws.Range("E38").Value = "  +0.87%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5921"
$ws.Range("E39").Value = "  -0.45%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.71"
$ws.Range("E40").Value = "  +0.92%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.904"
$ws.Range("E41").Value = "  -0.20%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1845"
$ws.Range("E42").Value = "  -0.48%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.441"
$ws.Range("E43").Value = "  +0.50%  "

# Row 44
$ws.Range("E44").Value = "  +0.29%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.38"
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07474"
$ws.Range("E46").Value = "  -3.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5556"
$ws.Range("E47").Value = "  -0.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.936"
$ws.Range("E48").Value = "  -0.87%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.11"
$ws.Range("E49").Value = "  +2.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.418"
$ws.Range("E50").Value = "  +3.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.99"
$ws.Range("E51").Value = "  -0.97%  "
